$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $escaped = $value -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range('D2') '29.714.47'
Set-TextValue $ws.Range('E2') '  +1.80%  '
Set-TextValue $ws.Range('D3') '1.856.82'
Set-TextValue $ws.Range('E3') '  +1.41%  '
Set-TextValue $ws.Range('D4') '0.9989'
Set-TextValue $ws.Range('E4') '  -0.02%  '
Set-TextValue $ws.Range('D5') '244.64'
Set-TextValue $ws.Range('E5') '  +0.67%  '
Set-TextValue $ws.Range('D6') '0.6406'
Set-TextValue $ws.Range('E6') '  +3.61%  '
Set-TextValue $ws.Range('D7') '0.9996'
Set-TextValue $ws.Range('E7') '  -0.06%  '
Set-TextValue $ws.Range('D8') '47.37'
Set-TextValue $ws.Range('E8') '  +4.55%  '
Set-TextValue $ws.Range('D9') '0.07517'
Set-TextValue $ws.Range('E9') '  +2.27%  '
Set-TextValue $ws.Range('D10') '0.2976'
Set-TextValue $ws.Range('E10') '  +2.33%  '
Set-TextValue $ws.Range('D11') '24.45'
Set-TextValue $ws.Range('E11') '  +5.31%  '
Set-TextValue $ws.Range('D12') '0.07662'
Set-TextValue $ws.Range('E12') '  -0.01%  '
Set-TextValue $ws.Range('D13') '1.870.41'
Set-TextValue $ws.Range('E13') '  +2.00%  '
Set-TextValue $ws.Range('D14') '5.035'
Set-TextValue $ws.Range('E14') '  +1.29%  '
Set-TextValue $ws.Range('D15') '0.6908'
Set-TextValue $ws.Range('E15') '  +3.48%  '
Set-TextValue $ws.Range('D16') '83.83'
Set-TextValue $ws.Range('E16') '  +1.67%  '
Set-TextValue $ws.Range('D17') '0.000009824'
Set-TextValue $ws.Range('E17') '  +9.93%  '
Set-TextValue $ws.Range('D18') '6.085'
Set-TextValue $ws.Range('E18') '  +4.14%  '
Set-TextValue $ws.Range('D19') '29.720.63'
Set-TextValue $ws.Range('E19') '  +1.86%  '
Set-TextValue $ws.Range('D20') '2.115.04'
Set-TextValue $ws.Range('E20') '  +1.38%  '
Set-TextValue $ws.Range('D21') '236.22'
Set-TextValue $ws.Range('E21') '  +0.61%  '
Set-TextValue $ws.Range('D22') '12.66'
Set-TextValue $ws.Range('E22') '  +1.56%  '
Set-TextValue $ws.Range('D23') '0.9997'
Set-TextValue $ws.Range('E23') '  -0.02%  '
Set-TextValue $ws.Range('D24') '7.496'
Set-TextValue $ws.Range('E24') '  +2.06%  '
Set-TextValue $ws.Range('D25') '0.9994'
Set-TextValue $ws.Range('E25') '  -0.10%  '
Set-TextValue $ws.Range('D26') '158.79'
Set-TextValue $ws.Range('E26') '  +0.44%  '
Set-TextValue $ws.Range('D27') '0.1419'
Set-TextValue $ws.Range('E27') '  +1.73%  '
Set-TextValue $ws.Range('D28') '8.526'
Set-TextValue $ws.Range('E28') '  -0.18%  '
Set-TextValue $ws.Range('D29') '17.90'
Set-TextValue $ws.Range('E29') '  +1.64%  '
Set-TextValue $ws.Range('D30') '0.06205'
Set-TextValue $ws.Range('E30') '  +6.68%  '
Set-TextValue $ws.Range('D31') '1.493'
Set-TextValue $ws.Range('E31') '  +0.26%  '
Set-TextValue $ws.Range('D32') '1.284'
Set-TextValue $ws.Range('E32') '  +6.15%  '
Set-TextValue $ws.Range('D33') '4.163'
Set-TextValue $ws.Range('E33') '  +1.97%  '
Set-TextValue $ws.Range('D34') '4.100'
Set-TextValue $ws.Range('E34') '  -0.02%  '
Set-TextValue $ws.Range('D35') '1.900'
Set-TextValue $ws.Range('E35') '  +2.77%  '
Set-TextValue $ws.Range('D36') '1.172'
Set-TextValue $ws.Range('E36') '  +2.83%  '
Set-TextValue $ws.Range('D37') '0.7286'
Set-TextValue $ws.Range('E37') '  +0.20%  '
Set-TextValue $ws.Range('D38') '2.602'
Set-TextValue $ws.Range('E38') '  -0.22%  '
Set-TextValue $ws.Range('D39') '2.830'
Set-TextValue $ws.Range('E39') '  -1.02%  '
Set-TextValue $ws.Range('D40') '0.01785'
Set-TextValue $ws.Range('E40') '  +1.87%  '
Set-TextValue $ws.Range('D41') '1.201.18'
Set-TextValue $ws.Range('E41') '  -1.60%  '
Set-TextValue $ws.Range('D42') '0.9222'
Set-TextValue $ws.Range('E42') '  +2.05%  '
Set-TextValue $ws.Range('D43') '6.232'
Set-TextValue $ws.Range('E43') '  +0.00%  '
Set-TextValue $ws.Range('D44') '0.9998'
Set-TextValue $ws.Range('E44') '  -0.04%  '
Set-TextValue $ws.Range('D45') '2.025.10'
Set-TextValue $ws.Range('E45') '  +1.80%  '
Set-TextValue $ws.Range('E46') '  +0.26%  '
Set-TextValue $ws.Range('D47') '66.54'
Set-TextValue $ws.Range('E47') '  +1.77%  '
Set-TextValue $ws.Range('D48') '0.00000000118'
Set-TextValue $ws.Range('E48') '  +0.65%  '
Set-TextValue $ws.Range('B49') 'EnergySwap'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D49') '9.223'
Set-TextValue $ws.Range('E49') '  +1.19%  '
Set-TextValue $ws.Range('B50') 'TheSandbox'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D50') '0.4059'
Set-TextValue $ws.Range('E50') '  +0.91%  '
Set-TextValue $ws.Range('D51') '0.05797'
Set-TextValue $ws.Range('E51') '  +0.94%  '
